$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1155.55
$ws.Range("C3").Value = 1666.23
$ws.Range("C4").Value = 2244.31
$ws.Range("C5").Value = 2871.13
